# Strike through the SP-section correction items that were addressed:
# "a.", "b.", "c." and "e." under "5.  SP:" get <w:strike/> applied to the
# paragraph mark and to every run in the paragraph (i.e. the whole
# paragraph, including its end-of-paragraph mark, is struck through).
# Items "d.", "f.", "g." and "h." are left untouched.

$d = $word.ActiveDocument

# Distinctive substrings that uniquely identify each target paragraph.
$targets = @(
    "No voy a eliminar cualquier usuario",
    "No se realiza validaci",
    "Para modificar un empleado",
    "Las transacciones se usan cuando tengo dos"
)

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    foreach ($t in $targets) {
        if ($text -like "*$t*") {
            $para.Range.Font.StrikeThrough = $true
            break
        }
    }
}
